# Insert a new data row before the existing row 66 ("Jengibre" weekly price
# records for Vega Central Mapocho de Santiago). Excel shifts the existing
# rows 66-122 down to 67-123 and carries their formatting (including the
# date-formatted column D style) along with them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new weekly record.
$ws.Range("A66").Value2 = 9
$ws.Range("B66").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C66").Value2 = "Metropolitana"
$ws.Range("D66").Value2 = 44977
$ws.Range("E66").Value2 = 13
$ws.Range("F66").Value2 = 100114007
$ws.Range("G66").Value2 = "Jengibre"
$ws.Range("H66").Value2 = "Sin especificar"
$ws.Range("I66").Value2 = "Primera"
$ws.Range("J66").Value2 = 520
$ws.Range("K66").Value2 = 18000
$ws.Range("L66").Value2 = 20000
$ws.Range("M66").Value2 = 19000
$ws.Range("N66").Value2 = "$/caja 13 kilos"
$ws.Range("O66").Value2 = "Perú"
$ws.Range("P66").Value2 = 1462
$ws.Range("Q66").Value2 = 13
$ws.Range("R66").Value2 = "Hortaliza"
